# Updated Master Gantt Chart
# Shorten the team-member name labels (drop last names) across every sheet
# that lists them: "Management Summary", "Gantt", "Meetings", "SA".

$wb = $excel.ActiveWorkbook

$targets = @(
    @{ Sheet = "Management Summary"; Cells = @{ A3="Jacob:"; A4="Cameron:"; A5="Conrad:"; A6="Benjamin:"; A7="Delaney:"; A8="Corbin:" }; Active = "A9" },
    @{ Sheet = "Gantt";              Cells = @{ A2="Jacob:"; A14="Cameron:"; A26="Conrad:"; A38="Benjamin:"; A50="Delaney:"; A62="Corbin:" }; Active = "A62" },
    @{ Sheet = "Meetings";           Cells = @{ A4="Jacob:"; A5="Cameron:"; A6="Conrad:"; A7="Benjamin:"; A8="Delaney:"; A9="Corbin:" }; Active = "A10" },
    @{ Sheet = "SA";                 Cells = @{ A2="Jacob:"; A6="Cameron:"; A10="Conrad:"; A14="Benjamin:"; A18="Delaney:"; A22="Corbin:" }; Active = "B29" }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    foreach ($addr in $t.Cells.Keys) {
        $ws.Range($addr).Value2 = $t.Cells[$addr]
    }
    # Touch each sheet so its own remembered cursor position updates too --
    # Excel keeps a per-sheet selection independent of which tab is active.
    $ws.Activate()
    $ws.Range($t.Active).Select()
}

# Finish on the "SA" sheet, which becomes the active tab after the edit.
$saSheet = $wb.Worksheets.Item("SA")
$saSheet.Activate()
$saSheet.Range("B29").Select()
